$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# 1. Rename the existing "表结构设计" sheet to "表结构设计-user"
$ws.Name = "表结构设计-user"

# 2. Add a new worksheet right after it, named "表结构设计－order"
$new = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws)
$new.Name = "表结构设计－order"

# 3. Fill in the new order-table sheet content
$new.Range("A2").Value = "客户订单表－order"
$new.Range("A3").Value = "序号"
$new.Range("B3").Value = "字段名"
$new.Range("C3").Value = "字段类型"
$new.Range("D3").Value = "字段说明"
$new.Range("A4").Value = 1
$new.Range("B4").Value = "orderid"
$new.Range("C4").Value = "varchar(32)"
$new.Range("D4").Value = "订单ID，主键，使用UUID"
$new.Range("A5").Value = 2
$new.Range("B5").Value = "linkman_name"
$new.Range("C5").Value = "varchar(32)"
$new.Range("D5").Value = "联系人姓名，默认是下单人的姓名，但是可以需改，不能为空"
$new.Range("A6").Value = 3
$new.Range("B6").Value = "phone"
$new.Range("C6").Value = "varchar(32)"
$new.Range("D6").Value = "联系电话，默认是下单人的手机号，不可为空"
$new.Range("A7").Value = 4
$new.Range("B7").Value = "show_addr"
$new.Range("C7").Value = "varchar(128)"
$new.Range("D7").Value = "演出地址"
$new.Range("A8").Value = 5
$new.Range("B8").Value = "show_date"
$new.Range("C8").Value = "datetime"
$new.Range("D8").Value = "演出时间，不可为空"
$new.Range("A9").Value = 6
$new.Range("B9").Value = "show_fee"
$new.Range("C9").Value = "int"
$new.Range("D9").Value = "预计演出金额"
$new.Range("A10").Value = 7
$new.Range("B10").Value = "request_mark"
$new.Range("C10").Value = "text"
$new.Range("D10").Value = "演出需求，由用户填写，需要什么类型的节目，看中哪些艺人"
$new.Range("A11").Value = 8
$new.Range("B11").Value = "userid"
$new.Range("C11").Value = "varchar(32)"
$new.Range("D11").Value = "下单人用户ID"
$new.Range("A12").Value = 9
$new.Range("B12").Value = "openid"
$new.Range("C12").Value = "varchar(45)"
$new.Range("D12").Value = "微信公众号关注用户的id"
$new.Range("A13").Value = 10
$new.Range("B13").Value = "create_date"
$new.Range("C13").Value = "datetime"
$new.Range("D13").Value = "定单创建时间"
$new.Range("A14").Value = 11
$new.Range("B14").Value = "modify_date"
$new.Range("C14").Value = "datetime"
$new.Range("D14").Value = "订单修改时间，客户代表谈完后修改金额之类的事宜"
$new.Range("A15").Value = 12
$new.Range("B15").Value = "modify_userid"
$new.Range("C15").Value = "varchar(32)"
$new.Range("D15").Value = "订单修改人ID"

# 4. Merge the title row
$new.Range("A2:D2").Merge()

# 5. Apply formatting to the new sheet, reusing existing formats so that
#    the style table grows in the same way the source workbook's did.
# Title bar (copy format from the existing "用户角色表" title on the user sheet)
$ws.Range("A18:D18").Copy()
$new.Range("A2:D2").PasteSpecial(-4122)

# Header row + index-number column (copy the centered/bordered style)
$ws.Range("A6").Copy()
$new.Range("A3:D3").PasteSpecial(-4122)
$new.Range("A4:A15").PasteSpecial(-4122)

# Data columns B:D (copy the bordered style, then left/vcenter align it)
$ws.Range("B6:D6").Copy()
$new.Range("B4:D15").PasteSpecial(-4122)
$new.Range("B4:D15").HorizontalAlignment = -4131

# 6. Column widths for the new sheet
$new.Columns("B:C").ColumnWidth = 13.5
$new.Columns("D").ColumnWidth = 55.5

# 7. Update the user-sheet header/index-number cells to use the new centered (no-vcenter) style
$ws.Range("A5:D5").HorizontalAlignment = -4108
$ws.Range("F5:I5").HorizontalAlignment = -4108
$ws.Range("F6:F18").HorizontalAlignment = -4108
$ws.Range("A19:D19").HorizontalAlignment = -4108

# 8. View/selection updates
$ws.Range("F29:F34").Select()
$new.Range("B8").Select()
$new.Activate()
